$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (46061 -> 46062) for every data row (rows 2 through 428).
for ($row = 2; $row -le 428; $row++) {
    $ws.Cells.Item($row, 3).Value = 46062
}
